$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string rich text, flattened) ---
$ws.Range("A8").Value = "Volume 30   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/24/2023  Through  4/30/2023"

# --- Cells changing from numeric to text "0" (shared string idx 20), style matches C14 ---
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("C14").Copy($ws.Range("C16"))
$ws.Range("C14").Copy($ws.Range("D17"))
$ws.Range("C14").Copy($ws.Range("D20"))
$ws.Range("C14").Copy($ws.Range("C23"))
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("C14").Copy($ws.Range("G27"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("C14").Copy($ws.Range("D29"))
$ws.Range("C14").Copy($ws.Range("G30"))

# --- Cells changing from numeric to text "***.*" (shared string idx 21), style matches E14 ---
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("E14").Copy($ws.Range("E17"))
$ws.Range("E14").Copy($ws.Range("E20"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("E14").Copy($ws.Range("H27"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("E14").Copy($ws.Range("E29"))
$ws.Range("E14").Copy($ws.Range("H30"))

# --- Cell changing from text "0" to numeric, style matches F23 (s=15) ---
$ws.Range("F23").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 2

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("N14").Value = -94.736842105263
$ws.Range("G15").Value = 3
$ws.Range("M15").Value = -88.888888888888
$ws.Range("N15").Value = -95.652173913043
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -100
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -33.333333333333
$ws.Range("J16").Value = 33
$ws.Range("K16").Value = -6.060606060606
$ws.Range("L16").Value = -11.428571428571
$ws.Range("M16").Value = -52.307692307692
$ws.Range("N16").Value = -87.398373983739
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -57.142857142857
$ws.Range("I17").Value = 61
$ws.Range("K17").Value = -22.784810126582
$ws.Range("L17").Value = 10.909090909090
$ws.Range("M17").Value = 38.636363636363
$ws.Range("N17").Value = -73.245614035087
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 4
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 26
$ws.Range("J18").Value = 31
$ws.Range("K18").Value = -16.129032258064
$ws.Range("L18").Value = -50.943396226415
$ws.Range("M18").Value = -7.142857142857
$ws.Range("N18").Value = -90.076335877862
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -28.571428571428
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = -38.709677419354
$ws.Range("I19").Value = 91
$ws.Range("J19").Value = 119
$ws.Range("K19").Value = -23.529411764705
$ws.Range("L19").Value = -2.150537634408
$ws.Range("M19").Value = 89.583333333333
$ws.Range("N19").Value = -18.75
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 35
$ws.Range("K20").Value = 66.666666666666
$ws.Range("L20").Value = 191.666666666667
$ws.Range("M20").Value = 66.666666666666
$ws.Range("N20").Value = -63.541666666666
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 11
$ws.Range("E21").Value = 9.090909090909
$ws.Range("F21").Value = 46
$ws.Range("G21").Value = 76
$ws.Range("H21").Value = -39.473684210526
$ws.Range("I21").Value = 246
$ws.Range("J21").Value = 290
$ws.Range("K21").Value = -15.172413793103
$ws.Range("L21").Value = -1.6
$ws.Range("M21").Value = 13.364055299539
$ws.Range("N21").Value = -75.050709939148
$ws.Range("D22").Value = 3
$ws.Range("G22").Value = 4
$ws.Range("J22").Value = 10
$ws.Range("K22").Value = 0
$ws.Range("C24").Value = 7
$ws.Range("E24").Value = -76.666666666666
$ws.Range("F24").Value = 47
$ws.Range("G24").Value = 108
$ws.Range("H24").Value = -56.481481481481
$ws.Range("I24").Value = 223
$ws.Range("J24").Value = 338
$ws.Range("K24").Value = -34.023668639053
$ws.Range("L24").Value = 18.617021276595
$ws.Range("M24").Value = 116.504854368932
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -75
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = -20
$ws.Range("I25").Value = 98
$ws.Range("J25").Value = 129
$ws.Range("K25").Value = -24.031007751938
$ws.Range("L25").Value = 27.272727272727
$ws.Range("M25").Value = -30.985915492957
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -25
$ws.Range("I26").Value = 6
$ws.Range("K26").Value = -45.454545454545
$ws.Range("L26").Value = 20
$ws.Range("F27").Value = 5
$ws.Range("I27").Value = 10
$ws.Range("K27").Value = -41.176470588235
$ws.Range("L27").Value = 25
$ws.Range("L28").Value = -50
$ws.Range("N28").Value = -94.736842105263
$ws.Range("L29").Value = -50
$ws.Range("N29").Value = -94.594594594594
